$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("M6").Value = "Jessica S. Tisch"
$ws.Range("A8").Value = "Volume 31   Number  48"
$ws.Range("C9").Value = "Report Covering the Week  11/25/2024  Through  12/1/2024"

# --- Crime Complaints data table (rows 14-30) ---
$ws.Range("N14").Value = -78.571428571428
$ws.Range("C15").Value = 1
$ws.Range("I15").Value = 38
$ws.Range("K15").Value = 80.95238095238
$ws.Range("L15").Value = 137.5
$ws.Range("M15").Value = 90
$ws.Range("N15").Value = 8.571428571428
$ws.Range("C16").Value = 4
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -20
$ws.Range("F16").Value = 20
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 33.333333333333
$ws.Range("I16").Value = 232
$ws.Range("J16").Value = 259
$ws.Range("K16").Value = -10.42471042471
$ws.Range("L16").Value = -6.072874493927
$ws.Range("M16").Value = -20.547945205479
$ws.Range("N16").Value = -62.998405103668
$ws.Range("C17").Value = 9
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = 125
$ws.Range("F17").Value = 32
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 430
$ws.Range("J17").Value = 385
$ws.Range("K17").Value = 11.688311688311
$ws.Range("L17").Value = 25.730994152046
$ws.Range("M17").Value = 76.954732510288
$ws.Range("N17").Value = 33.956386292834
$ws.Range("C18").Value = 7
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = 16.666666666666
$ws.Range("F18").Value = 14
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 16.666666666666
$ws.Range("I18").Value = 153
$ws.Range("J18").Value = 218
$ws.Range("K18").Value = -29.816513761467
$ws.Range("L18").Value = 26.446280991735
$ws.Range("M18").Value = -57.734806629834
$ws.Range("N18").Value = -88.91304347826
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = -15.384615384615
$ws.Range("F19").Value = 50
$ws.Range("G19").Value = 68
$ws.Range("H19").Value = -26.470588235294
$ws.Range("I19").Value = 781
$ws.Range("J19").Value = 635
$ws.Range("K19").Value = 22.992125984252
$ws.Range("L19").Value = 36.300174520069
$ws.Range("M19").Value = 88.192771084337
$ws.Range("N19").Value = 43.830570902394
$ws.Range("D20").Value = 12
$ws.Range("E20").Value = -8.333333333333
$ws.Range("F20").Value = 37
$ws.Range("H20").Value = -15.90909090909
$ws.Range("I20").Value = 424
$ws.Range("J20").Value = 488
$ws.Range("K20").Value = -13.11475409836
$ws.Range("L20").Value = 28.875379939209
$ws.Range("M20").Value = 101.904761904762
$ws.Range("N20").Value = -76.085730400451
$ws.Range("C21").Value = 43
$ws.Range("D21").Value = 40
$ws.Range("E21").Value = 7.5
$ws.Range("F21").Value = 159
$ws.Range("G21").Value = 171
$ws.Range("H21").Value = -7.017543859649
$ws.Range("I21").Value = 2061
$ws.Range("J21").Value = 2011
$ws.Range("K21").Value = 2.486325211337
$ws.Range("L21").Value = 26.055045871559
$ws.Range("M21").Value = 32.796391752577
$ws.Range("N21").Value = -56.083528659705
$ws.Range("C22").NumberFormat = '#,##0'
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 2
$ws.Range("I22").Value = 23
$ws.Range("K22").Value = 109.090909090909
$ws.Range("L22").Value = 15
$ws.Range("M22").Value = 27.777777777777
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = -14.285714285714
$ws.Range("I23").Value = 99
$ws.Range("J23").Value = 110
$ws.Range("K23").Value = -10
$ws.Range("L23").Value = -1
$ws.Range("M23").Value = 67.796610169491
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -25.806451612903
$ws.Range("F24").Value = 91
$ws.Range("G24").Value = 106
$ws.Range("H24").Value = -14.150943396226
$ws.Range("I24").Value = 1159
$ws.Range("J24").Value = 1479
$ws.Range("K24").Value = -21.636240703177
$ws.Range("L24").Value = -6.305578011317
$ws.Range("M24").Value = 34.298957126303
$ws.Range("C25").Value = 6
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -57.142857142857
$ws.Range("F25").Value = 28
$ws.Range("G25").Value = 43
$ws.Range("H25").Value = -34.883720930232
$ws.Range("I25").Value = 436
$ws.Range("J25").Value = 635
$ws.Range("K25").Value = -31.338582677165
$ws.Range("L25").Value = -11.382113821138
$ws.Range("C26").Value = 15
$ws.Range("E26").Value = 66.666666666666
$ws.Range("F26").Value = 51
$ws.Range("G26").Value = 44
$ws.Range("H26").Value = 15.90909090909
$ws.Range("I26").Value = 531
$ws.Range("J26").Value = 491
$ws.Range("K26").Value = 8.146639511201
$ws.Range("L26").Value = 9.484536082474
$ws.Range("M26").Value = -11.647254575707
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = "'0"
$ws.Range("E27").Value = "'***.*"
$ws.Range("I27").Value = 44
$ws.Range("K27").Value = 29.411764705882
$ws.Range("L27").Value = 22.222222222222
$ws.Range("C28").Value = 1
$ws.Range("D28").NumberFormat = '#,##0'
$ws.Range("D28").Value = 1
$ws.Range("E28").NumberFormat = '#,##0.0;"-"#,##0.0'
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 8
$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 100
$ws.Range("I28").Value = 65
$ws.Range("J28").Value = 63
$ws.Range("K28").Value = 3.174603174603
$ws.Range("L28").Value = 30
$ws.Range("D29").Value = "'0"
$ws.Range("E29").Value = "'***.*"
$ws.Range("N29").Value = -66.666666666666
$ws.Range("D30").Value = "'0"
$ws.Range("E30").Value = "'***.*"
$ws.Range("N30").Value = -65.90909090909
